# related work about digital twin
# Adds rows describing additional digital-twin-related references to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 3 : white paper summary (D3 = year, E3 = abstract, justified text)
# ---------------------------------------------------------------------
$ws.Range("D3").Value = 2020
$ws.Range("E3").Value = "本白皮书通过梳理数字孪生技术和产业应用发展情况，分析数字孪生技术热点、行业动态和未来趋势，提出相关的标准化工作需求，希望可以作为数字孪生技术领域、产业发展和标准化之间的初始连接纽带，加快推动数字孪生发展应用。"
$ws.Range("E3").Font.Name = "等线"
$ws.Range("E3").Font.Size = 10.5
$ws.Range("E3").HorizontalAlignment = -4130
$ws.Range("E3").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Row 4 : Alibaba IDC demo (E4 = mixed-font title, F4 = link)
# ---------------------------------------------------------------------
$ws.Range("E4").Value = "阿里巴巴的IDC数字孪生"
$ws.Range("E4").HorizontalAlignment = -4108
$ws.Range("E4").VerticalAlignment = -4108
$ws.Range("E4").Characters(6, 3).Font.Name = "Times New Roman"
$ws.Range("E4").Characters(9, 4).Font.Name = "宋体"
$ws.Range("F4").Value = "https://market.m.taobao.com/app/txddp/idc-digital-twin-pr/index.html?spm=a2c6h.12873639.0.0.78714addNSn6j5#/park/EA133/D"

# ---------------------------------------------------------------------
# Row 5 : oneNET platform demo (E5 = mixed-font title, F5 = link)
# ---------------------------------------------------------------------
$ws.Range("E5").Value = "oneNET平台的数据孪生可视化demo"
$ws.Range("E5").Characters(7, 10).Font.Name = "宋体"
$ws.Range("E5").Characters(17, 4).Font.Name = "Times New Roman"
$ws.Range("F5").Value = "https://open.iot.10086.cn/view/main/#/share2d?id=5dc3c5fdf03f1e16d7b8975d&shareCode=1234"

# ---------------------------------------------------------------------
# Row 6 : Mumiantree platform demo (E6 = mixed-font title, F6 = link)
# ---------------------------------------------------------------------
$ws.Range("E6").Value = "国内木棉树平台的数字孪生demo"
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E6").VerticalAlignment = -4108
$ws.Range("E6").Characters(13, 4).Font.Name = "Times New Roman"
$ws.Range("F6").Value = "http://www.mms3d.cn/html/product/slxm/"

# ---------------------------------------------------------------------
# Row 7 : aircraft structural life prediction paper (D7 = year, E7 = abstract)
# ---------------------------------------------------------------------
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = "Reengineering of the aircraft structural life prediction process to fully exploit advances in very high performance digital computing is proposed. The proposed process utilizes an ultrahigh fidelity model of individual aircraft by tail number, a Digital Twin, to integrate computation of structural deflections and temperatures in response to flight conditions, with resulting local damage and material state evolution. A conceptual model of how the Digital Twin can be used for predicting the life of aircraft structure and assuring its structural integrity is presented. The technical challenges to developing and deploying a Digital Twin are discussed in detail."
$ws.Range("E7").Font.Name = "等线"
$ws.Range("E7").Font.Size = 10.5
$ws.Range("E7").HorizontalAlignment = -4130
$ws.Range("E7").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Row 8 : big data and digital twin in smart manufacturing paper
# ---------------------------------------------------------------------
$ws.Range("D8").Value = 2017
$ws.Range("E8").Value = "With the advances in new-generation information technologies, especially big data and digital twin, smart manufacturing is becoming the focus of global manufacturing transformation and upgrading. Intelligence comes from data. Integrated analysis for the manufacturing big data is beneficial to all aspects of manufacturing. Besides, the digital twin paves a way for the cyber-physical integration of manufacturing, which is an important bottleneck to achieve smart manufacturing. In this paper, the big data and digital twin in manufacturing are reviewed, including their concept as well as their applications in product design, production planning, manufacturing, and predictive maintenance. On this basis, the similarities and differences between big data and digital twin are compared from the general and data perspectives. Since the big data and digital twin can be complementary, how they can be integrated to promote smart manufacturing are discussed."
$ws.Range("E8").Font.Name = "等线"
$ws.Range("E8").Font.Size = 10.5
$ws.Range("E8").HorizontalAlignment = -4130
$ws.Range("E8").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Row 9 : experimentable digital twins paper (E9 = abstract, F9 = link)
# ---------------------------------------------------------------------
$ws.Range("D9").Value = 2018
$ws.Range("E9").Value = "Digital twins represent real objects or subjects with their data, functions, and communication capabilities in the digital world. As nodes within the internet of things, they enable networking and thus the automation of complex value-added chains. The application of simulation techniques brings digital twins to life and makes them experimentable; digital twins become experimentable digital twins (EDTs). Initially, these EDTs communicate with each other purely in the virtual world. The resulting networks of interacting EDTs model different application scenarios and are simulated in virtual testbeds, providing new foundations for comprehensive simulation-based systems engineering. Its focus is on EDTs, which become more detailed with every single application. Thus, complete digital representations of the respective real assets and their behaviors are created successively. The networking of EDTs with real assets leads to hybrid application scenarios in which EDTs are used in combination with real hardware, thus realizing complex control algorithms, innovative user interfaces, or mental models for intelligent systems."
$ws.Range("E9").Font.Name = "Arial"
$ws.Range("E9").Font.Size = 11.5
$ws.Range("E9").Font.Color = 3355443
$ws.Range("E9").HorizontalAlignment = -4130
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("F9").Value = "https://ieeexplore.ieee.org/document/8289327"

# Row heights to match the wrapped/justified content of the new rows.
$ws.Rows.Item(3).RowHeight = 27.6
$ws.Rows.Item(6).RowHeight = 14.4
$ws.Rows.Item(7).RowHeight = 55.2
$ws.Rows.Item(8).RowHeight = 69
$ws.Rows.Item(9).RowHeight = 100.8

$ws.Range("F11").Select()
